$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 rows above the existing table to make room for a small
# "cover sheet" block (Company Name / Your Name / Phone Number / Email ID).
$ws.Rows("1:4").Insert()

# Keep Table1 anchored on the (now shifted) header/data rows.
$ws.ListObjects("Table1").Resize($ws.Range("A5:O6"))

# Row 1 - Company Name:
$r = $ws.Range("A1")
$r.Interior.Color = 65535
$r.Font.Bold = $true
$r.Value = "Company Name:"
$ws.Range("B1").Interior.ThemeColor = 6

# Row 2 - Your Name:
$r = $ws.Range("A2")
$r.Interior.ThemeColor = 6
$r.Font.Bold = $true
$r.Value = "Your Name:"
$ws.Range("B2").Interior.Color = 65535

# Row 3 - Phone Number:
$r = $ws.Range("A3")
$r.Interior.Color = 65535
$r.Font.Bold = $true
$r.Value = "Phone Number:"
$ws.Range("B3").Interior.ThemeColor = 6

# Row 4 - Email ID:
$r = $ws.Range("A4")
$r.Interior.ThemeColor = 6
$r.Font.Bold = $true
$r.Value = "Email ID:"
$ws.Range("B4").Interior.Color = 65535

$ws.Range("B12").Select()
